$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)   # "Previously added"
$ws2 = $wb.Worksheets.Item(2)   # "New"

$xlPasteFormats = -4122

# ---------------------------------------------------------------------------
# Step 1: snapshot the 5 rows currently sitting in "New" (rows 2..6) - these
# are the rows that are being archived into "Previously added".
# ---------------------------------------------------------------------------
$firstOldRow = 2
$lastOldRow = 6

$archiveRows = New-Object System.Collections.ArrayList
for ($r = $firstOldRow; $r -le $lastOldRow; $r++) {
    $row = @{
        A = $ws2.Cells.Item($r, 1).Value2
        B = $ws2.Cells.Item($r, 2).Value2
        C = $ws2.Cells.Item($r, 3).Value2
        D = $ws2.Cells.Item($r, 4).Value2
        E = $ws2.Cells.Item($r, 5).Value2
        F = $ws2.Cells.Item($r, 6).Value2
    }
    $archiveRows.Add($row) | Out-Null
}

# ---------------------------------------------------------------------------
# Step 2: append those rows to the end of "Previously added".
# ---------------------------------------------------------------------------
$destFirstRow = $ws1.UsedRange.Rows.Count + 1   # 267

for ($i = 0; $i -lt $archiveRows.Count; $i++) {
    $destRow = $destFirstRow + $i
    $data = $archiveRows[$i]

    # Bring over formatting from the row directly above so styles (s="3"/"4"/"2")
    # match the rest of the table exactly.
    $ws1.Range("A" + ($destRow - 1) + ":F" + ($destRow - 1)).Copy()
    $ws1.Range("A" + $destRow + ":F" + $destRow).PasteSpecial($xlPasteFormats)

    # Force column E to text first in case the cadastre number looks numeric.
    $ws1.Range("E" + $destRow).NumberFormat = "@"

    $ws1.Cells.Item($destRow, 1).Value = $data.A
    $ws1.Cells.Item($destRow, 2).Value = $data.B
    $ws1.Cells.Item($destRow, 3).Value = $data.C
    $ws1.Cells.Item($destRow, 4).Value = $data.D
    $ws1.Cells.Item($destRow, 5).Value = $data.E
    $ws1.Cells.Item($destRow, 6).Value = $data.F

    $ws1.Hyperlinks.Add($ws1.Range("A" + $destRow), $data.A) | Out-Null

    # Re-apply the row formatting once more: Hyperlinks.Add() stamps its own
    # built-in "Hyperlink" style onto column A, and the NumberFormat tweak
    # above can nudge the style index too, so restore everything from the
    # same template row as the final step.
    $ws1.Range("A" + ($destRow - 1) + ":F" + ($destRow - 1)).Copy()
    $ws1.Range("A" + $destRow + ":F" + $destRow).PasteSpecial($xlPasteFormats)
}

# ---------------------------------------------------------------------------
# Step 3: remove the archived rows from "New", dropping the now-stale
# hyperlinks along with them.
# ---------------------------------------------------------------------------
$ws2.Hyperlinks.Delete()
$ws2.Range("A" + ($firstOldRow + 1) + ":F" + $lastOldRow).EntireRow.Delete()

# ---------------------------------------------------------------------------
# Step 4: write the single brand-new listing into "New" row 2.
# ---------------------------------------------------------------------------
$newRow = 2
$newLink = "https://www.ss.com/msg/lv/real-estate/wood/talsi-and-reg/dundagas-pag/mjdlj.html"
$newPrice = "169 000 €"
$newDistrict = "Talsi un raj."
$newArea = "59 ha."
$newCadastre = "88500020198"
$newDate = 45982.368055555555

$ws2.Range("E" + $newRow).NumberFormat = "@"

$ws2.Cells.Item($newRow, 1).Value = $newLink
$ws2.Cells.Item($newRow, 2).Value = $newPrice
$ws2.Cells.Item($newRow, 3).Value = $newDistrict
$ws2.Cells.Item($newRow, 4).Value = $newArea
$ws2.Cells.Item($newRow, 5).Value = $newCadastre
$ws2.Cells.Item($newRow, 6).Value = $newDate

$ws2.Hyperlinks.Add($ws2.Range("A" + $newRow), $newLink) | Out-Null

# Restore the table's normal row styling (Hyperlinks.Add recolors column A).
$ws1.Range("A266:F266").Copy()
$ws2.Range("A" + $newRow + ":F" + $newRow).PasteSpecial($xlPasteFormats)

$excel.CutCopyMode = $false
